# TeamANSReport.xlsx — "Add files via upload"
#
# The underlying data edit captured by the diff is a Story re-assignment:
# the Backlog/Sprint1 row that used to point at US11 "No bigamy" now points
# at US10 "Marriage after 14", and its status flips from "Coding" to "Done"
# (shared-string indices 123->122, 77->76, 198->192 in the OOXML package).
# Everything else in the xml_diff (xr:* namespace/uid churn, fileVersion,
# calcId, window/selection geometry, default row/col sizing) is Excel-build
# resave noise, not a content change, so only the cell values below are
# touched here.

$wb = $excel.ActiveWorkbook

# --- Sprint1 sheet ---------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")

# Row 10: Story ID + Story Name move from US11/"No bigamy" to
# US10/"Marriage after 14".
$sprint1.Range("A10").Value = "US10"
$sprint1.Range("B10").Value = "Marriage after 14"

# Rows 9 and 10: Status moves from "Coding" to "Done".
$sprint1.Range("D9").Value = "Done"
$sprint1.Range("D10").Value = "Done"

# --- Backlog sheet -----------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")

# Row 9: Story ID + Story Name move from US11/"No bigamy" to
# US10/"Marriage after 14" (mirrors the Sprint1 change above).
$backlog.Range("B9").Value = "US10"
$backlog.Range("C9").Value = "Marriage after 14"
